$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared-string text blocks introduced by this edit.
$signinText = @'
type: signin
width: 2
height: 1
h3: Raise / Sponsor Funds
p: Try out our new feature. Raise Funds for your next project or Join us in distributing rations.
button.primary: Create a Ticket*goto("/createticket")
button.secondary: View Tickets*goto("/tickets")
svg: /icons/stars.svg
'@

$blogText = @'
type: blog
width: 2
height: 1
ser: 144
'@

# Build new row 10, reusing existing card text where the value repeats.
$ws.Range("A10").Value = 43971
$ws.Range("A10").NumberFormat = $ws.Range("A9").NumberFormat

$ws.Range("C10").Value = $signinText
$ws.Range("B10").Value = $blogText
$ws.Range("D10").Value = $ws.Range("E9").Value2
$ws.Range("E10").Value = $ws.Range("F9").Value2
$ws.Range("F10").Value = $ws.Range("G9").Value2
$ws.Range("G10").Value = $ws.Range("C9").Value2
$ws.Range("H10").Value = $ws.Range("I9").Value2
$ws.Range("I10").Value = $ws.Range("H9").Value2
$ws.Range("J10").Value = $ws.Range("K9").Value2
$ws.Range("K10").Value = $ws.Range("M9").Value2

$ws.Range("B10:K10").WrapText = $true
$ws.Rows.Item(10).RowHeight = 255

# Update the view to match the post-edit selection/scroll state.
$ws.Range("K10").Select()
$excel.ActiveWindow.ScrollColumn = 2
